$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert new record at row 10 (most recent week for Carahue) ---
$ws.Rows.Item(10).Insert()

$ws.Cells.Item(10,1).Value = 9
$ws.Cells.Item(10,2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(10,3).Value = "Metropolitana"
$ws.Cells.Item(10,4).Value = 44552
$ws.Cells.Item(10,5).Value = 13
$ws.Cells.Item(10,6).Value = 100112022
$ws.Cells.Item(10,7).Value = "Arveja Verde"
$ws.Cells.Item(10,8).Value = "Sin especificar"
$ws.Cells.Item(10,9).Value = "Primera"
$ws.Cells.Item(10,10).Value = 52
$ws.Cells.Item(10,11).Value = 12000
$ws.Cells.Item(10,12).Value = 13000
$ws.Cells.Item(10,13).Value = 12500
$ws.Cells.Item(10,14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(10,15).Value = "Carahue"
$ws.Cells.Item(10,16).Value = 500
$ws.Cells.Item(10,17).Value = 25
$ws.Cells.Item(10,18).Value = "Hortaliza"

# --- Insert new record at row 18 (new weekly entry, Región Metropolitana) ---
$ws.Rows.Item(18).Insert()

$ws.Cells.Item(18,1).Value = 9
$ws.Cells.Item(18,2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(18,3).Value = "Metropolitana"
$ws.Cells.Item(18,4).Value = 44497
$ws.Cells.Item(18,5).Value = 13
$ws.Cells.Item(18,6).Value = 100112022
$ws.Cells.Item(18,7).Value = "Arveja Verde"
$ws.Cells.Item(18,8).Value = "Sin especificar"
$ws.Cells.Item(18,9).Value = "Primera"
$ws.Cells.Item(18,10).Value = 35
$ws.Cells.Item(18,11).Value = 1300
$ws.Cells.Item(18,12).Value = 1500
$ws.Cells.Item(18,13).Value = 1414
$ws.Cells.Item(18,14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(18,15).Value = "Región Metropolitana"
$ws.Cells.Item(18,16).Value = 57
$ws.Cells.Item(18,17).Value = 25
$ws.Cells.Item(18,18).Value = "Hortaliza"
